$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap data rows 2 and 3 (keep header row 1 and column A trial index unchanged)
$row2 = @($ws.Range("B2").Value2, $ws.Range("C2").Value2, $ws.Range("D2").Value2, $ws.Range("E2").Value2, $ws.Range("F2").Value2, $ws.Range("G2").Value2, $ws.Range("H2").Value2, $ws.Range("I2").Value2, $ws.Range("J2").Value2)
$row3 = @($ws.Range("B3").Value2, $ws.Range("C3").Value2, $ws.Range("D3").Value2, $ws.Range("E3").Value2, $ws.Range("F3").Value2, $ws.Range("G3").Value2, $ws.Range("H3").Value2, $ws.Range("I3").Value2, $ws.Range("J3").Value2)

$ws.Range("B2").Value = $row3[0]
$ws.Range("C2").Value = $row3[1]
$ws.Range("D2").Value = $row3[2]
$ws.Range("E2").Value = $row3[3]
$ws.Range("F2").Value = $row3[4]
$ws.Range("G2").Value = $row3[5]
$ws.Range("H2").Value = $row3[6]
$ws.Range("I2").Value = $row3[7]
$ws.Range("J2").Value = $row3[8]

$ws.Range("B3").Value = $row2[0]
$ws.Range("C3").Value = $row2[1]
$ws.Range("D3").Value = $row2[2]
$ws.Range("E3").Value = $row2[3]
$ws.Range("F3").Value = $row2[4]
$ws.Range("G3").Value = $row2[5]
$ws.Range("H3").Value = $row2[6]
$ws.Range("I3").Value = $row2[7]
$ws.Range("J3").Value = $row2[8]
